$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q (year 2020) ---

# Q3: empty cell, same formatting as the other row-3 border cells (e.g. A3 / D3:P3)
$ws.Range("A3").Copy()
$ws.Range("Q3").PasteSpecial(-4122) # xlPasteFormats

# Q4: header value 2020, formatted like the other year headers (D4:P4) but
# with top-aligned text instead of center-aligned.
$ws.Range("D4").Copy()
$ws.Range("Q4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Q4").VerticalAlignment = -4160 # xlTop
$ws.Range("Q4").Value = 2020

# Q5: data value for the first indicator row, keep the existing style (s=13)
$ws.Range("Q5").Value = 1.1000000000000001

# Q6: data value for the second indicator row, formatted with one decimal
# place like the other numeric data cells, based on the row-3 border style.
$ws.Range("A3").Copy()
$ws.Range("Q6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Q6").NumberFormat = "0.0"
$ws.Range("Q6").Value = 7

# --- Selection / active cell ---
$ws.Range("J22").Select()
